$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "culture_collection" was re-added to this MIxS plant-associated template in
# an earlier revision; per INSDC2017 review it is being removed again. That
# means deleting its whole column (header row 15 holds the field names) so
# everything to the right shifts one column to the left.

$headerRow = 15
$targetName = "culture_collection"

# Locate the last used column and the column to delete dynamically, so the
# script keeps working even if the layout shifts slightly.
$lastCol = $ws.UsedRange.Columns.Count
$delCol = -1
for ($c = 1; $c -le $lastCol; $c++) {
  if ($ws.Cells.Item($headerRow, $c).Text -eq $targetName) {
    $delCol = $c
  }
}

if ($delCol -eq -1) {
  throw "Could not find the '$targetName' column to delete."
}

# Cell comments are anchored to fixed addresses, so Range.Delete() does not
# carry them along with the shifting cells the way values/styles do. Capture
# every comment from the doomed column through the end of the used range
# before deleting, so they can be re-applied one column to the left.
$texts = @{}
for ($c = $delCol; $c -le $lastCol; $c++) {
  $cell = $ws.Cells.Item($headerRow, $c)
  $cm = $cell.Comment
  if ($cm -ne $null) {
    $texts[$c] = $cm.Text()
  }
}

# Delete the whole column; values, styles and shared strings shift left on
# their own.
$ws.Cells.Item($headerRow, $delCol).EntireColumn.Delete()

# Re-home the captured comments: whatever used to sit in column c+1 now
# lives in column c. Clear any leftover comment first, since the delete
# above left the originals exactly where they were.
for ($c = $delCol; $c -le ($lastCol - 1); $c++) {
  $cell = $ws.Cells.Item($headerRow, $c)
  if ($cell.Comment -ne $null) {
    $cell.Comment.Delete()
  }
  $srcCol = $c + 1
  if ($texts.ContainsKey($srcCol)) {
    $cell.AddComment($texts[$srcCol])
  }
}

# The very last (now out-of-range) column may still carry its old comment;
# drop it since that column no longer exists in the shifted layout.
$strayCell = $ws.Cells.Item($headerRow, $lastCol)
if ($strayCell.Comment -ne $null) {
  $strayCell.Comment.Delete()
}
